# financial-statement-analysis workbook touch-up
#  - rename "Balance sheet" -> "Balance Sheet"
#  - fix the "Shareholders' Equity" label (drop the apostrophe) on the Balance Sheet
#  - move the active tab from Cash Flow Statement to Balance Sheet, and
#    update each sheet's remembered selection accordingly

$wb = $excel.ActiveWorkbook

$incomeStatement = $wb.Worksheets.Item(1)
$balanceSheet    = $wb.Worksheets.Item(2)
$cashFlow        = $wb.Worksheets.Item(3)

# Rename the Balance sheet tab (capitalisation fix)
$balanceSheet.Name = "Balance Sheet"

# Correct the shared-string label used in A6 of the Balance Sheet
$balanceSheet.Range("A6").Value = "Shareholders Equity"

# Cash Flow Statement keeps its existing selection (C5); only the
# Balance Sheet's remembered selection moves, to A6
$balanceSheet.Range("A6").Select()

# Make the Balance Sheet the active/selected tab instead of Cash Flow Statement
$balanceSheet.Activate()
